# lut_truth_table.xlsx edit:
# Set the E-column "select" inputs to 1 for rows 6-13 (s0..s3 pattern 0000..0111)
# and clear the three test markers (rows 16 and 21) that previously held the
# extra "fizz"/"buzz"/"fizz-buzz" notes in column F.
# This changes the active SUM-mode test pattern from {E11,E16,E21} to {E6..E13}.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows whose E value should now be 1 (select bit on)
$onRows = 6,7,8,9,10,11,12,13
foreach ($r in $onRows) {
    $ws.Cells.Item($r, 5).Value = 1
}

# Rows whose E value should now be cleared (select bit off)
$offRows = 16,21
foreach ($r in $offRows) {
    $ws.Cells.Item($r, 5).ClearContents()
}

# Column F previously carried "fizz" / "buzz" / "fizz-buzz" note text on
# rows 9, 11, 15, 16, 18, 21 - clear all of them.
$noteRows = 9,11,15,16,18,21
foreach ($r in $noteRows) {
    $ws.Cells.Item($r, 6).ClearContents()
}

# Update the saved selection to E14 (matches the diff's sheetView selection)
$ws.Range("E14").Select() | Out-Null

$wb.Save()
